$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8170576691627502
$ws.Range("B1").Value = 1.347966194152832
$ws.Range("C1").Value = 3.590045213699341
$ws.Range("D1").Value = 2.644761562347412
$ws.Range("E1").Value = 1.657705903053284
